$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.315.02"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.911.79"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'0.723"
$ws.Range("E5").Value = "  +9.42%  "
$ws.Range("D6").Value = "'253.68"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'40.66"
$ws.Range("E8").Value = "  -1.45%  "
$ws.Range("E9").Value = "  +3.76%  "
$ws.Range("D10").Value = "'52.27"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").Value = "'0.0769"
$ws.Range("E11").Value = "  +7.16%  "
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").Value = "2.185.47"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "'12.77"
$ws.Range("E14").Value = "  +5.41%  "
$ws.Range("D15").Value = "'0.721"
$ws.Range("E15").Value = "  +2.58%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'4.93"
$ws.Range("E16").Value = "  +1.37%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.888.22"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "35.265.79"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "'74.39"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "0.0₃0850"
$ws.Range("E20").Value = "  +2.37%  "
$ws.Range("D21").Value = "'243.90"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("D22").Value = "'13.08"
$ws.Range("E22").Value = "  +4.56%  "
$ws.Range("D23").Value = "'5.08"
$ws.Range("E23").Value = "  +4.90%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "'2.46"
$ws.Range("E25").Value = "  +4.48%  "
$ws.Range("D26").Value = "'2.38"
$ws.Range("E26").Value = "  +3.75%  "
$ws.Range("D27").Value = "'166.67"
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").Value = "'8.68"
$ws.Range("E28").Value = "  +2.01%  "
$ws.Range("D29").Value = "'18.76"
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("D30").Value = "'0.133"
$ws.Range("E30").Value = "  +4.60%  "
$ws.Range("D31").Value = "4.129.85"
$ws.Range("E31").Value = "  +19.49%  "
$ws.Range("D32").Value = "'4.35"
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("E33").Value = "  +14.48%  "
$ws.Range("D34").Value = "'1.64"
$ws.Range("E34").Value = "  +23.17%  "
$ws.Range("D35").Value = "'0.0584"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("D36").Value = "'4.21"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("E37").Value = "  +0.14%  "
$ws.Range("D38").Value = "'0.921"
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").Value = "'2.04"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'17.19"
$ws.Range("E40").Value = "  +5.11%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0216"
$ws.Range("E41").Value = "  +3.71%  "
$ws.Range("D42").Value = "'96.67"
$ws.Range("E42").Value = "  +7.23%  "
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("E44").Value = "  -0.97%  "
$ws.Range("D45").Value = "1.338.49"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'2.42"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "'6.77"
$ws.Range("E48").Value = "  +2.81%  "
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "'45.44"
$ws.Range("E50").Value = "  -5.53%  "
$ws.Range("D51").Value = "'12.01"
$ws.Range("E51").Value = "  +15.34%  "
